$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1452
$wsExhibit.Range("F5").Value = 14

# Sheet "全部类型" (All types) - mirrors the same data, apply identical update
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1452
$wsAll.Range("F5").Value = 14
